# Update "想去人数" (want-to-go count) figures in the 杭州-漫展信息 workbook.
# Sheet "展览" (Exhibitions) and its rollup sheet "全部类型" (All types) share
# the same events, and sheet "本地生活" (Local life) also has one matching
# event that needs the same refresh.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsLocal   = $wb.Worksheets.Item("本地生活")
$wsAll     = $wb.Worksheets.Item("全部类型")

# 展览 sheet updates
$wsExhibit.Range("F7").Value  = 970
$wsExhibit.Range("F18").Value = 2577
$wsExhibit.Range("F20").Value = 1073
$wsExhibit.Range("F21").Value = 3558
$wsExhibit.Range("F23").Value = 837
$wsExhibit.Range("F29").Value = 361
$wsExhibit.Range("F30").Value = 195
$wsExhibit.Range("F32").Value = 1325
$wsExhibit.Range("F33").Value = 1936
$wsExhibit.Range("F35").Value = 32

# 本地生活 sheet update
$wsLocal.Range("F2").Value = 417

# 全部类型 sheet updates
$wsAll.Range("F2").Value  = 417
$wsAll.Range("F6").Value  = 970
$wsAll.Range("F21").Value = 2577
$wsAll.Range("F23").Value = 1073
$wsAll.Range("F24").Value = 3558
$wsAll.Range("F26").Value = 837
$wsAll.Range("F37").Value = 361
$wsAll.Range("F38").Value = 195
$wsAll.Range("F41").Value = 1325
$wsAll.Range("F42").Value = 1936
